$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2 through 17 changes from 45243 (2023-11-13)
# to 45244 (2023-11-14), keeping the existing date formatting/style intact.
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
